# Add Santander Parser row + update chunk-splitting table range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 10) --------------------------------------------------
$ws.Range("A10").Value = "0069-100128-251"
$ws.Range("B10").Value = "Santander"
$ws.Range("C10").Value = "SINGAPORE"
$ws.Range("D10").Value = "(SG)"
$ws.Range("E10").Value = "C.P.Trading"
$ws.Range("F10").Value = "Saving"
$ws.Range("G10").Value = "US"

# --- Grow Table1 so the new row participates in the table/autofilter -------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G10"))

# --- Highlight the new bank's branch cell (bigger font, box border, centered)
$c10 = $ws.Range("C10")
$c10.Font.Name = "Calibri"
$c10.Font.Size = 14
$c10.HorizontalAlignment = -4108
$c10.Borders.LineStyle = 1
$c10.Borders.Weight = 2

# --- Row 10 is taller to fit the bigger font --------------------------------
$ws.Rows.Item(10).RowHeight = 18.15

# --- Leave the selection on the newly edited cell ---------------------------
$ws.Range("C10").Select()
